# Applies the "getting closer on postgame hitter" edit to the Loveless, Ethan
# postgame hitter report. Updates play-by-play details (pitcher name, pitch
# mix, fastball velo, counts, exit velo / launch angle, hit type/result) for
# several at-bats, and clears a handful of exit-velo / launch-angle / hit-type
# cells back to "not yet computed" (empty inlineStr) placeholders.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- At-bat block starting row 10 (Inning 3) ---
$ws.Range("M10").Value = ""

# --- At-bat block starting row 12 ---
$ws.Range("M12").Value = ""

# --- At-bat block starting row 17 ---
$ws.Range("J17").Value = "FB,CB,CH"

# --- At-bat block starting row 19 ---
$ws.Range("J19").Value = 4
$ws.Range("M19").Value = ""
$ws.Range("J20").Value = 0
$ws.Range("M21").Value = ""

# --- At-bat block starting row 23 ---
$ws.Range("J23").Value = "Roblez"
$ws.Range("M23").Value = ""
$ws.Range("J25").Value = "88-90 MPH"
$ws.Range("J26").Value = "FB,CB,CH"

# --- At-bat block starting row 28 ---
$ws.Range("J28").Value = 5
$ws.Range("M28").Value = ""
$ws.Range("J29").Value = 2
$ws.Range("M30").Value = ""

# --- At-bat block starting row 32 ---
$ws.Range("J32").Value = "Herbst"
$ws.Range("M32").Value = "Ground Ball"
$ws.Range("M33").Value = "Double"
$ws.Range("J34").Value = "83-85 MPH"
$ws.Range("J35").Value = "SL,FB,CB,CH"

# --- At-bat block starting row 37 ---
$ws.Range("J37").Value = 7
$ws.Range("M37").Value = ""
$ws.Range("J38").Value = 0
$ws.Range("M39").Value = ""

# --- At-bat block starting row 41 ---
$ws.Range("J41").Value = "Plum"
$ws.Range("M41").Value = ""
$ws.Range("M42").Value = "Undefined"
$ws.Range("J43").Value = "84-86 MPH"
$ws.Range("J44").Value = "SL,FB,CH"

# --- At-bat block starting row 46 ---
$ws.Range("J46").Value = 7
$ws.Range("M46").Value = ""
$ws.Range("J47").Value = 2
$ws.Range("M48").Value = ""

# --- At-bat block starting row 50 ---
$ws.Range("J50").Value = "Plum"
$ws.Range("M50").Value = ""
$ws.Range("J51").Value = "Right"
$ws.Range("M51").Value = "Undefined"
$ws.Range("J52").Value = "84-86 MPH"

# --- At-bat block starting row 61 ---
$ws.Range("J61").Value = 9
$ws.Range("M61").Value = ""
$ws.Range("J62").Value = 1
$ws.Range("M63").Value = ""

# --- At-bat block starting row 65 ---
$ws.Range("J65").Value = "Thompson"
$ws.Range("M65").Value = "Popup"
$ws.Range("J66").Value = "Left"
$ws.Range("M66").Value = "Out"
$ws.Range("J67").Value = "84-84 MPH"
$ws.Range("J68").Value = "SL,FB,CH"
